$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Text content updates (progress notes / version tags) for the GCA
# Improvement List, per the 1.5.5 status update commit.
# ---------------------------------------------------------------------------

# Row 2 - Thorough testing of permutations
$ws.Range("E2").Value = "1.5.5*"
$ws.Range("F2").Value = "Finished testing the various trace combinations, but I couldn't reproduce the intermittent error when plotting vs Pout. I did fix an issue you saw when moving from one X-axis type to another (say Pin to Frequency, and vice versa). Unfortunately the fix for now is for me to delete then recreate the trace entirely. I expect these trace changes will only need to be done during the initial configuration, so hopefully this isn't a big issue."

# Row 3 - Increase the number of points in power interpolation plots
$ws.Range("E3").Value = "1.5.5"
$ws.Range("F3").Value = "vs Pin, Pout plots have 10x points"

# Row 4 - End of sweep conditions
$ws.Range("F4").Value = "This depends on the development team in Munich. I can make the request, but I don't know if/when it would get implemented. A work-around Greg recommended is to display the power range in the diagram title. Not sure if that would work since in general the ranges in the diagram may not all be the same. Another idea is to use markers at the trace endpoints."

# Row 5 - Move plots to a different figure
$ws.Range("E5").Value = "1.5.5"

# Row 6 - Allow plot of gain and power values on the same figure
$ws.Range("E6").Value = "1.5.5"

# Row 7 - Resolve the (x-axis) scaling issue
$ws.Range("B7").Value = "Resolve the x-axis scaling issue caused by a change in measurement conditions"
$ws.Range("F7").Value = "I was able to reproduce and document this phenomenon. Unfortunately, it is an issue with our VNA firmware. I generated a report to send to our firmware developers. I will update you when I get a response. A work-around would be for me to delete then recreate the trace entirely, but this would mean that you'd have to re-setup the plots (not convenient)."

# Row 8 - Message when plot of value exceeding the measurement boundaries
$ws.Range("D8").Value = "Not started"
$ws.Range("E8").Value = "1.5.6"
$ws.Range("F8").Value = "The GUI will not let you input incorrect values: any at value that is not possible is rounded to the nearest possible value. If you watch the value as you enter it, you should see this. That said, I will add a warning in the next version."

# Row 9 - Change plot nomenclature
$ws.Range("E9").Value = "1.5.5"

# Row 10 - CSV file comma-delimeter
$ws.Range("E10").Value = "1.5.5"

# Row 11 - CSV file add settings as header
$ws.Range("E11").Value = "1.5.6"

# Row 12 - CSV file add pulse definition info
$ws.Range("E12").Value = "1.5.6"

# Row 13 - CSV file add phase of Pout
$ws.Range("E13").Value = "1.5.5"
$ws.Range("F13").Value = "FWIW, also added phase of Pin (0 deg)"

# ---------------------------------------------------------------------------
# Column E ("Proposed Date for Resolution / SW Version") is now populated
# for most rows, so center it like the other status/category columns. E2
# previously stood out in bold; normalize it to match the rest of the column.
# ---------------------------------------------------------------------------
$ws.Range("E2").Font.Bold = $false
$ws.Range("E2:E17").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Row heights: a few rows grew to fit the updated / longer comments.
# ---------------------------------------------------------------------------
$ws.Rows(3).RowHeight = 19.899999999999999
$ws.Rows(4).RowHeight = 90
$ws.Rows(7).RowHeight = 90
$ws.Rows(8).RowHeight = 60

# ---------------------------------------------------------------------------
# Selection moved to E7 (editor was last looking at the updated comment).
# ---------------------------------------------------------------------------
$ws.Range("E7").Select() | Out-Null
